# Add a new "NET" ticker row to both the IncomeReport and EarningsReport
# sheets (row 4), mirroring the existing CRWD / RUN rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: IncomeReport
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("IncomeReport")

# Duplicate the formatting (fill/border/alignment) of row 2 (style "2",
# the same style alternation NET should get as the 3rd data row) onto the
# new row 4 before writing values into it.
$ws1.Range("A2:I2").Copy()
$ws1.Range("A4:I4").PasteSpecial(-4122)
$ws1.Rows.Item(4).RowHeight = 30

$ws1.Cells.Item(4, 1).Value = "NET"
$ws1.Cells.Item(4, 2).Value = "212.17M -> 234.52M -> 253.86M -> 274.7M -> 290.18M"
$ws1.Cells.Item(4, 3).Value = 37
$ws1.Cells.Item(4, 4).Value = "-0.13 -> -0.2 -> -0.13 -> -0.14 -> -0.12"
$ws1.Cells.Item(4, 5).Value = 8
$ws1.Cells.Item(4, 6).Value = "(54.95M) -> 6.65M -> 6.13M -> 34.08M -> 20.81M"
$ws1.Cells.Item(4, 7).Value = 138
$ws1.Cells.Item(4, 8).Value = "5.75 <- N/A <- N/A <- N/A"
$ws1.Cells.Item(4, 9).Value = "20.63 <- 16.39 <- 21.85 <- 18.96"

# Column D widened slightly (20.29 -> ~20.71 characters).
$ws1.Columns.Item(4).ColumnWidth = 19.8

# ---------------------------------------------------------------------
# Sheet 2: EarningsReport
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("EarningsReport")

$ws2.Range("A2:G2").Copy()
$ws2.Range("A4:G4").PasteSpecial(-4122)
$ws2.Rows.Item(4).RowHeight = 30

$ws2.Cells.Item(4, 1).Value = "NET"
$ws2.Cells.Item(4, 2).Value = 0.03
$ws2.Cells.Item(4, 3).Value = 0.08
$ws2.Cells.Item(4, 4).Value = 166
$ws2.Cells.Item(4, 5).Value = "158, 157"
$ws2.Cells.Item(4, 6).Value = "30, 31"
$ws2.Cells.Item(4, 7).Value = 59

Write-Output "NET row appended to IncomeReport and EarningsReport"
